{"js": "const replacements = [\n  [\"2024-10-07 Monday\", \"2024-10-08 Tuesday\"],\n  [\"257\u00d72=514\", \"728\u00d79=6552\"],\n  [\"154\u00d75=770\", \"953\u00d74=3812\"],\n  [\"978\u00d75=4890\", \"809\u00d79=7281\"],\n  [\"844\u00d79=7596\", \"795\u00d72=1590\"],\n  [\"966\u00d72=1932\", \"684\u00d72=1368\"],\n  [\"170\u00d77=1190\", \"956\u00d76=5736\"],\n  [\"854\u00d77=5978\", \"590\u00d77=4130\"],\n  [\"259\u00d74=1036\", \"255\u00d73=765\"],\n  [\"614\u00d73=1842\", \"405\u00d74=1620\"],\n  [\"167\u00d73=501\", \"941\u00d77=6587\"],\n  [\"760\u00d77=5320\", \"984\u00d75=4920\"],\n  [\"380\u00d79=3420\", \"168\u00d77=1176\"],\n  [\"936\u00d74=3744\", \"224\u00d74=896\"],\n  [\"820\u00d73=2460\", \"985\u00d73=2955\"],\n  [\"390\u00d72=780\", \"365\u00d75=1825\"],\n  [\"236\u00d78=1888\", \"346\u00d73=1038\"],\n  [\"533\u00d73=1599\", \"645\u00d72=1290\"],\n  [\"724\u00d72=1448\", \"648\u00d73=1944\"],\n  [\"181\u00d77=1267\", \"525\u00d74=2100\"],\n  [\"924\u00d79=8316\", \"765\u00d74=3060\"],\n  [\"797\u00d74=3188\", \"401\u00d79=3609\"],\n  [\"128\u00d77=896\", \"303\u00d76=1818\"],\n  [\"865\u00d79=7785\", \"906\u00d74=3624\"],\n  [\"106\u00d74=424\", \"218\u00d76=1308\"],\n  [\"807\u00d76=4842\", \"225\u00d73=675\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-07 Monday\", \"2024-10-08 Tuesday\"),\n    @(\"257\u00d72=514\", \"728\u00d79=6552\"),\n    @(\"154\u00d75=770\", \"953\u00d74=3812\"),\n    @(\"978\u00d75=4890\", \"809\u00d79=7281\"),\n    @(\"844\u00d79=7596\", \"795\u00d72=1590\"),\n    @(\"966\u00d72=1932\", \"684\u00d72=1368\"),\n    @(\"170\u00d77=1190\", \"956\u00d76=5736\"),\n    @(\"854\u00d77=5978\", \"590\u00d77=4130\"),\n    @(\"259\u00d74=1036\", \"255\u00d73=765\"),\n    @(\"614\u00d73=1842\", \"405\u00d74=1620\"),\n    @(\"167\u00d73=501\", \"941\u00d77=6587\"),\n    @(\"760\u00d77=5320\", \"984\u00d75=4920\"),\n    @(\"380\u00d79=3420\", \"168\u00d77=1176\"),\n    @(\"936\u00d74=3744\", \"224\u00d74=896\"),\n    @(\"820\u00d73=2460\", \"985\u00d73=2955\"),\n    @(\"390\u00d72=780\", \"365\u00d75=1825\"),\n    @(\"236\u00d78=1888\", \"346\u00d73=1038\"),\n    @(\"533\u00d73=1599\", \"645\u00d72=1290\"),\n    @(\"724\u00d72=1448\", \"648\u00d73=1944\"),\n    @(\"181\u00d77=1267\", \"525\u00d74=2100\"),\n    @(\"924\u00d79=8316\", \"765\u00d74=3060\"),\n    @(\"797\u00d74=3188\", \"401\u00d79=3609\"),\n    @(\"128\u00d77=896\", \"303\u00d76=1818\"),\n    @(\"865\u00d79=7785\", \"906\u00d74=3624\"),\n    @(\"106\u00d74=424\", \"218\u00d76=1308\"),\n    @(\"807\u00d76=4842\", \"225\u00d73=675\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n$d.Save()\n"}
